$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.998.36"
$ws.Cells.Item(2, 5).Value = "  -0.36%  "

$ws.Cells.Item(3, 4).Value = "1.741.74"
$ws.Cells.Item(3, 5).Value = "  +1.23%  "

$ws.Cells.Item(4, 4).Value = "'1.001"
$ws.Cells.Item(4, 4).NumberFormat = "General"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.47%  "

$ws.Cells.Item(5, 4).Value = "'311.56"
$ws.Cells.Item(5, 4).NumberFormat = "General"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -2.24%  "

$ws.Cells.Item(6, 5).Value = "  -0.49%  "

$ws.Cells.Item(7, 4).Value = "'0.4988"
$ws.Cells.Item(7, 4).NumberFormat = "General"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +6.83%  "

$ws.Cells.Item(8, 4).Value = "'0.3596"
$ws.Cells.Item(8, 4).NumberFormat = "General"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +4.28%  "

$ws.Cells.Item(9, 4).Value = "'42.57"
$ws.Cells.Item(9, 4).NumberFormat = "General"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +0.51%  "

$ws.Cells.Item(10, 4).Value = "'0.07270"
$ws.Cells.Item(10, 4).NumberFormat = "General"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -0.56%  "

$ws.Cells.Item(11, 4).Value = "'1.062"
$ws.Cells.Item(11, 4).NumberFormat = "General"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.74%  "

$ws.Cells.Item(12, 4).Value = "'1.001"
$ws.Cells.Item(12, 4).NumberFormat = "General"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -0.40%  "

$ws.Cells.Item(13, 4).Value = "'20.27"
$ws.Cells.Item(13, 4).NumberFormat = "General"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +1.73%  "

$ws.Cells.Item(14, 4).Value = "'5.969"
$ws.Cells.Item(14, 4).NumberFormat = "General"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +1.48%  "

$ws.Cells.Item(15, 4).Value = "1.739.01"
$ws.Cells.Item(15, 5).Value = "  +0.09%  "

$ws.Cells.Item(16, 4).Value = "'6.859"
$ws.Cells.Item(16, 4).NumberFormat = "General"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.73%  "

$ws.Cells.Item(17, 4).Value = "'86.78"
$ws.Cells.Item(17, 4).NumberFormat = "General"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -3.51%  "

$ws.Cells.Item(18, 4).Value = "'0.00001037"
$ws.Cells.Item(18, 4).NumberFormat = "General"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -1.10%  "

$ws.Cells.Item(19, 4).Value = "'0.06369"
$ws.Cells.Item(19, 4).NumberFormat = "General"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.88%  "

$ws.Cells.Item(20, 5).Value = "  -0.69%  "

$ws.Cells.Item(21, 4).Value = "'16.59"
$ws.Cells.Item(21, 4).NumberFormat = "General"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +0.59%  "

$ws.Cells.Item(22, 5).Value = "  +1.55%  "

$ws.Cells.Item(23, 4).Value = "27.076.55"
$ws.Cells.Item(23, 5).Value = "  -0.26%  "

$ws.Cells.Item(24, 4).Value = "'11.29"
$ws.Cells.Item(24, 4).NumberFormat = "General"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +4.12%  "

$ws.Cells.Item(25, 4).Value = "'2.044"
$ws.Cells.Item(25, 4).NumberFormat = "General"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -3.85%  "

$ws.Cells.Item(26, 4).Value = "'155.06"
$ws.Cells.Item(26, 4).NumberFormat = "General"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -1.22%  "

$ws.Cells.Item(27, 4).Value = "'19.93"
$ws.Cells.Item(27, 4).NumberFormat = "General"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +1.93%  "

$ws.Cells.Item(28, 4).Value = "1.941.07"
$ws.Cells.Item(28, 5).Value = "  +0.46%  "

$ws.Cells.Item(29, 4).Value = "'2.153"
$ws.Cells.Item(29, 4).NumberFormat = "General"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.42%  "

$ws.Cells.Item(30, 4).Value = "'120.88"
$ws.Cells.Item(30, 4).NumberFormat = "General"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +1.10%  "

$ws.Cells.Item(31, 5).Value = "  +2.72%  "

$ws.Cells.Item(32, 4).Value = "'0.09588"
$ws.Cells.Item(32, 4).NumberFormat = "General"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +5.30%  "

$ws.Cells.Item(33, 4).Value = "'3.573"
$ws.Cells.Item(33, 4).NumberFormat = "General"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.85%  "

$ws.Cells.Item(34, 4).Value = "'5.392"
$ws.Cells.Item(34, 4).NumberFormat = "General"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +1.00%  "

$ws.Cells.Item(35, 4).Value = "'0.02203"
$ws.Cells.Item(35, 4).NumberFormat = "General"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.29%  "

$ws.Cells.Item(36, 4).Value = "'0.05870"
$ws.Cells.Item(36, 4).NumberFormat = "General"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +0.47%  "

$ws.Cells.Item(37, 4).Value = "'11.07"
$ws.Cells.Item(37, 4).NumberFormat = "General"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.72%  "

$ws.Cells.Item(38, 5).Value = "  +2.26%  "

$ws.Cells.Item(39, 4).Value = "'0.2002"
$ws.Cells.Item(39, 4).NumberFormat = "General"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +0.30%  "

$ws.Cells.Item(40, 4).Value = "'4.774"
$ws.Cells.Item(40, 4).NumberFormat = "General"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +0.03%  "

$ws.Cells.Item(41, 4).Value = "'0.6030"
$ws.Cells.Item(41, 4).NumberFormat = "General"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +0.99%  "

$ws.Cells.Item(42, 4).Value = "'1.110"
$ws.Cells.Item(42, 4).NumberFormat = "General"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -2.34%  "

$ws.Cells.Item(43, 4).Value = "'7.554"
$ws.Cells.Item(43, 4).NumberFormat = "General"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +0.54%  "

$ws.Cells.Item(44, 4).Value = "'12.94"
$ws.Cells.Item(44, 4).NumberFormat = "General"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +2.45%  "

$ws.Cells.Item(45, 4).Value = "'3.599"
$ws.Cells.Item(45, 4).NumberFormat = "General"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -1.42%  "

$ws.Cells.Item(46, 4).Value = "'0.5665"
$ws.Cells.Item(46, 4).NumberFormat = "General"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +1.01%  "

$ws.Cells.Item(47, 4).Value = "'120.36"
$ws.Cells.Item(47, 4).NumberFormat = "General"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.75%  "

$ws.Cells.Item(48, 4).Value = "'1.864"
$ws.Cells.Item(48, 4).NumberFormat = "General"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +0.03%  "

$ws.Cells.Item(49, 5).Value = "  +1.26%  "

$ws.Cells.Item(50, 4).Value = "'0.06677"
$ws.Cells.Item(50, 4).NumberFormat = "General"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.10%  "

$ws.Cells.Item(51, 5).Value = "  -0.55%  "
